# Bump the "Förändrad" (Changed) date in column C from 46074 to 46075
# for every data row (rows 2 through 246) on the active worksheet.
# Value2() is used (rather than Value()) because column C is
# date-formatted, and Value() would hand back a DateTime instead of the
# underlying serial-date number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2() -eq 46074) {
        $cell.Value2() = 46075
    }
}
